# Sync file from Google Drive
# Refresh the live "NextBus" feed snapshot (EstimatedTimeOfArrival in column F,
# MinutesToArrival in column O) across the three NextBus sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "NextBus1" ---
$ws = $wb.Worksheets.Item("NextBus1")
$ws.Cells.Item(2, 6).Value  = 45689.83564814815
$ws.Cells.Item(2, 15).Value = 7
$ws.Cells.Item(3, 6).Value  = 45689.84452546296
$ws.Cells.Item(3, 15).Value = 20
$ws.Cells.Item(4, 6).Value  = 45689.8408449074
$ws.Cells.Item(4, 15).Value = 14
$ws.Cells.Item(5, 6).Value  = 45689.83836805556
$ws.Cells.Item(5, 15).Value = 11
$ws.Cells.Item(6, 6).Value  = 45689.83313657407
$ws.Cells.Item(6, 15).Value = 3
$ws.Cells.Item(7, 6).Value  = 45689.82979166666
$ws.Cells.Item(7, 15).Value = -1
$ws.Cells.Item(8, 6).Value  = 45689.83210648148
$ws.Cells.Item(8, 15).Value = 2
$ws.Cells.Item(9, 6).Value  = 45689.83795138889
$ws.Cells.Item(9, 15).Value = 10
$ws.Cells.Item(10, 6).Value  = 45689.83380787037
$ws.Cells.Item(10, 15).Value = 4
$ws.Cells.Item(11, 6).Value  = 45689.83043981482
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(12, 6).Value  = 45689.83019675926
$ws.Cells.Item(12, 15).Value = 0
$ws.Cells.Item(13, 6).Value  = 45689.83017361111
$ws.Cells.Item(13, 15).Value = 0
$ws.Cells.Item(14, 6).Value  = 45689.83678240741
$ws.Cells.Item(14, 15).Value = 8
$ws.Cells.Item(15, 6).Value  = 45689.83122685185
$ws.Cells.Item(15, 15).Value = 0

# --- Sheet "NextBus2" ---
$ws = $wb.Worksheets.Item("NextBus2")
$ws.Cells.Item(2, 6).Value  = 45689.84717592593
$ws.Cells.Item(2, 15).Value = 23
$ws.Cells.Item(3, 15).Value = 35
$ws.Cells.Item(4, 6).Value  = 45689.85758101852
$ws.Cells.Item(4, 15).Value = 38
$ws.Cells.Item(5, 6).Value  = 45689.84912037037
$ws.Cells.Item(5, 15).Value = 26
$ws.Cells.Item(6, 6).Value  = 45689.84199074074
$ws.Cells.Item(6, 15).Value = 16
$ws.Cells.Item(7, 6).Value  = 45689.83717592592
$ws.Cells.Item(7, 15).Value = 9
$ws.Cells.Item(8, 6).Value  = 45689.83760416666
$ws.Cells.Item(8, 15).Value = 10
$ws.Cells.Item(9, 6).Value  = 45689.84760416667
$ws.Cells.Item(9, 15).Value = 24
$ws.Cells.Item(10, 6).Value  = 45689.8431712963
$ws.Cells.Item(10, 15).Value = 18
$ws.Cells.Item(11, 6).Value  = 45689.83763888889
$ws.Cells.Item(11, 15).Value = 10
$ws.Cells.Item(12, 6).Value  = 45689.84039351852
$ws.Cells.Item(13, 6).Value  = 45689.83671296296
$ws.Cells.Item(13, 15).Value = 8
$ws.Cells.Item(14, 6).Value  = 45689.84362268518
$ws.Cells.Item(14, 15).Value = 18
$ws.Cells.Item(15, 6).Value  = 45689.8390625
$ws.Cells.Item(15, 15).Value = 12

# --- Sheet "NextBus3" ---
$ws = $wb.Worksheets.Item("NextBus3")
$ws.Cells.Item(2, 15).Value = 33
$ws.Cells.Item(3, 15).Value = 41
$ws.Cells.Item(4, 6).Value  = 45689.85752314814
$ws.Cells.Item(4, 15).Value = 38
$ws.Cells.Item(5, 15).Value = 21
$ws.Cells.Item(6, 6).Value  = 45689.84678240741
$ws.Cells.Item(6, 15).Value = 23
$ws.Cells.Item(7, 6).Value  = 45689.86015046296
$ws.Cells.Item(7, 15).Value = 42
$ws.Cells.Item(8, 6).Value  = 45689.85244212963
$ws.Cells.Item(8, 15).Value = 31
$ws.Cells.Item(9, 15).Value = 24
$ws.Cells.Item(10, 15).Value = 22
$ws.Cells.Item(11, 6).Value  = 45689.84806712963
$ws.Cells.Item(11, 15).Value = 25
$ws.Cells.Item(12, 6).Value  = 45689.85505787037
$ws.Cells.Item(12, 15).Value = 35
$ws.Cells.Item(13, 6).Value  = 45689.8527662037
$ws.Cells.Item(13, 15).Value = 31
